$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 34,4
$data[0,0] = 45847.60416666666
$data[0,1] = 283.9400024414062
$data[0,2] = 283.885009765625
$data[0,3] = 282.9173840958459
$data[1,0] = 45847.64583333334
$data[1,1] = 283.2200012207031
$data[1,2] = 283.9400024414062
$data[1,3] = 294.485928141514
$data[2,0] = 45847.6875
$data[2,1] = 283.5799865722656
$data[2,2] = 283.2200012207031
$data[2,3] = 278.66844766402
$data[3,0] = 45847.72916666666
$data[3,1] = 283.0199890136719
$data[3,2] = 283.5799865722656
$data[3,3] = 278.6722002174581
$data[4,0] = 45847.77083333334
$data[4,1] = 282.7900085449219
$data[4,2] = 283.0199890136719
$data[4,3] = 279.846387524925
$data[5,0] = 45847.8125
$data[5,1] = 283.0400085449219
$data[5,2] = 282.7900085449219
$data[5,3] = 286.7017277114304
$data[6,0] = 45848.5625
$data[6,1] = 285.6199951171875
$data[6,2] = 283.0400085449219
$data[6,3] = 283.0819289931062
$data[7,0] = 45848.60416666666
$data[7,1] = 286.9500122070312
$data[7,2] = 285.6199951171875
$data[7,3] = 276.7136237333615
$data[8,0] = 45848.64583333334
$data[8,1] = 287.2099914550781
$data[8,2] = 286.9500122070312
$data[8,3] = 285.9455160428381
$data[9,0] = 45848.6875
$data[9,1] = 287.3500061035156
$data[9,2] = 287.2099914550781
$data[9,3] = 285.544055400027
$data[10,0] = 45848.72916666666
$data[10,1] = 287.5599975585938
$data[10,2] = 287.3500061035156
$data[10,3] = 288.041512361188
$data[11,0] = 45848.77083333334
$data[11,1] = 287.9450073242188
$data[11,2] = 287.5599975585938
$data[11,3] = 283.1223701249718
$data[12,0] = 45848.8125
$data[12,1] = 288.0950012207031
$data[12,2] = 287.9450073242188
$data[12,3] = 289.3544638487778
$data[13,0] = 45849.5625
$data[13,1] = 285.3599853515625
$data[13,2] = 288.0950012207031
$data[13,3] = 282.5655250833711
$data[14,0] = 45849.60416666666
$data[14,1] = 286.0950012207031
$data[14,2] = 285.3599853515625
$data[14,3] = 294.4982778093452
$data[15,0] = 45849.64583333334
$data[15,1] = 286.25
$data[15,2] = 286.0950012207031
$data[15,3] = 286.9545434440288
$data[16,0] = 45849.6875
$data[16,1] = 286.7000122070312
$data[16,2] = 286.25
$data[16,3] = 294.0256109447321
$data[17,0] = 45849.72916666666
$data[17,1] = 286.5450134277344
$data[17,2] = 286.7000122070312
$data[17,3] = 286.2124809091559
$data[18,0] = 45849.77083333334
$data[18,1] = 286.9505920410156
$data[18,2] = 286.5450134277344
$data[18,3] = 280.9820255135897
$data[19,0] = 45849.8125
$data[19,1] = 286.7099914550781
$data[19,2] = 286.9505920410156
$data[19,3] = 297.1051758949924
$data[20,0] = 45852.5625
$data[20,1] = 287.1499938964844
$data[20,2] = 286.7099914550781
$data[20,3] = 284.3581615373283
$data[21,0] = 45852.60416666666
$data[21,1] = 288.1199951171875
$data[21,2] = 287.1499938964844
$data[21,3] = 286.4229552677923
$data[22,0] = 45852.64583333334
$data[22,1] = 288.2699890136719
$data[22,2] = 288.1199951171875
$data[22,3] = 284.7697823027374
$data[23,0] = 45852.6875
$data[23,1] = 289.1099853515625
$data[23,2] = 288.2699890136719
$data[23,3] = 285.4067617662777
$data[24,0] = 45852.72916666666
$data[24,1] = 288.4400024414062
$data[24,2] = 289.1099853515625
$data[24,3] = 288.6322468223087
$data[25,0] = 45852.77083333334
$data[25,1] = 288.3999938964844
$data[25,2] = 288.4400024414062
$data[25,3] = 288.6033369807128
$data[26,0] = 45852.8125
$data[26,1] = 288.6799926757812
$data[26,2] = 288.3999938964844
$data[26,3] = 277.5393936133036
$data[27,0] = 45853.5625
$data[27,1] = 287.5799865722656
$data[27,2] = 288.6799926757812
$data[27,3] = 278.4177521205455
$data[28,0] = 45853.60416666666
$data[28,1] = 287.0499877929688
$data[28,2] = 287.5799865722656
$data[28,3] = 297.0202864193549
$data[29,0] = 45853.64583333334
$data[29,1] = 286.739990234375
$data[29,2] = 287.0499877929688
$data[29,3] = 291.6851385350887
$data[30,0] = 45853.6875
$data[30,1] = 286.4252014160156
$data[30,2] = 286.739990234375
$data[30,3] = 288.8684018728767
$data[31,0] = 45853.72916666666
$data[31,1] = 286.6000061035156
$data[31,2] = 286.4252014160156
$data[31,3] = 281.9299627668753
$data[32,0] = 45853.77083333334
$data[32,1] = 286.7049865722656
$data[32,2] = 286.6000061035156
$data[32,3] = 286.2776708586654
$data[33,0] = 45853.8125
$data[33,1] = 286.2999877929688
$data[33,2] = 286.7049865722656
$data[33,3] = 289.0007722774114

$ws.Range("A2:D35").Value = $data

$ws.Rows.Item(36).Delete()